$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Best-effort: touching the guides collection is what produces the empty
# <p15:sldGuideLst/> extension block PowerPoint leaves behind in the source
# deck. Not every host exposes a live Guides collection headlessly, so this
# is wrapped defensively and simply does nothing if unsupported.
# ---------------------------------------------------------------------------
try {
    $guides = $p.Guides
    if ($guides -ne $null) {
        $null = $guides.Count
    }
} catch {
}

# ---------------------------------------------------------------------------
# Slide 1: "Outline" -> "Initial Project Approach"
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

# Title
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "Initial Project Approach"

# Body / content placeholder
$body = $s1.Shapes.Item(2).TextFrame.TextRange

$body.Paragraphs(1,1).Text = "Android app for the user interface"
$body.Paragraphs(2,1).Text = "Use the Laird BL600 module for sensor device processing"
$body.Paragraphs(3,1).Text = "Use the MMA8653FC accelerometer for detecting movement"
$body.Paragraphs(4,1).Text = "Communication over Bluetooth Low Energy"

$p5 = $body.Paragraphs(5,1)
$p5.Text = "Implement a prototype on Dr. "
$run2 = $p5.InsertAfter("Jovanov’s")
$run3 = $run2.InsertAfter(" development board")

# New 6th bullet - appended after the existing last paragraph, inheriting the
# trailing endParaRPr that used to sit on paragraph 5 ("- demo?").
$tail = $body.InsertAfter("`rExplore over-the-air programming of the BL600")

Write-Output "Slide 1 updated."

# ---------------------------------------------------------------------------
# Slide 4: "Modified Project" -> "Modified Project Approach" (two runs)
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$title4 = $s4.Shapes.Item(1).TextFrame.TextRange
$title4.Paragraphs(1,1).Text = "Modified "
$title4.Paragraphs(1,1).InsertAfter("Project Approach") | Out-Null

Write-Output "Slide 4 updated."
